# Add a new log entry row (row 32) to the tracker sheet, matching the
# style/format of the previous row, and update the sheet selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 31
$newRow = 32

# New date entry (one day after the previous log entry: 2025-01-07 -> 2025-01-08)
$ws.Cells.Item($newRow, 1).Value = 45665
$ws.Cells.Item($newRow, 1).NumberFormat = $ws.Cells.Item($lastRow, 1).NumberFormat

# Task description and hours for the new entry
$ws.Cells.Item($newRow, 2).Value = "Adding more content to the current area "
$ws.Cells.Item($newRow, 3).Value = 6

# Update the sheet's active selection
$ws.Range("K22").Select()
